# Apply cryptos.xlsx data refresh (updated prices / 1h volume changes / row reorder for 47-50)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.309.60'
$ws.Range('E2').Value = '  -2.64%  '
$ws.Range('D3').Value = '3.683.63'
$ws.Range('E3').Value = '  -3.33%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '684.45'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.70%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '160.13'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -6.36%  '
$ws.Range('D7').Value = '3.681.65'
$ws.Range('E7').Value = '  -3.37%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.493'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -6.18%  '
$ws.Range('E10').Value = '  -8.82%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.21'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.30%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.436'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -10.43%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000232'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -7.07%  '
$ws.Range('D14').Value = '4.307.66'
$ws.Range('E14').Value = '  -3.38%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '32.41'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -11.29%  '
$ws.Range('D16').Value = '3.681.64'
$ws.Range('E16').Value = '  -3.18%  '
$ws.Range('D17').Value = '69.342.52'
$ws.Range('E17').Value = '  -2.87%  '
$ws.Range('E18').Value = '  -1.33%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '15.89'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -9.30%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.43'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -10.91%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '472.46'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -7.43%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.88'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.43%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.646'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -9.64%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '79.56'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -5.16%  '
$ws.Range('D25').Value = '3.830.95'
$ws.Range('E25').Value = '  -3.26%  '
$ws.Range('E26').Value = '  -0.10%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000124'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -11.92%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.87'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -14.54%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.14'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -11.18%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.70'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -10.64%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.74'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -13.08%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.02'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -10.58%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.64'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -9.88%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.00'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.30%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '26.62'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -8.56%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.159'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -7.44%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '8.17'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -12.28%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.06'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -8.85%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.25'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -6.03%  '
$ws.Range('E40').Value = '  -0.01%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0901'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -10.71%  '
$ws.Range('E42').Value = '  +0.01%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.940'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -7.07%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '164.92'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.85%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '47.89'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.15%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.71'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -15.53%  '
$ws.Range('B47').Value = 'ONDO'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.31'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -6.39%  '
$ws.Range('B48').Value = 'SuiNetwork'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.10'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.69%  '
$ws.Range('B49').Value = 'FLOKI'
$ws.Range('C49').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.000272'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -10.08%  '
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '27.87'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -8.56%  '
$ws.Range('E51').Value = '  -9.38%  '
